# "Generate Report for handback"
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the first data row
# (the row for the 53c9dd87-... file) on both the "zh-cn" and "de-de"
# worksheets, reflecting the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-08 14:28:38"
$wsZhCn.Range("G2").Value = "2016-01-08 14:29:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-08 14:28:52"
$wsDeDe.Range("G2").Value = "2016-01-08 14:29:52"
